$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell "D2" "61.630.17"
Set-TextCell "E2" "  -3.00%  "

Set-TextCell "D3" "2.575.23"
Set-TextCell "E3" "  -5.37%  "

Set-TextCell "E4" "  +0.03%  "

Set-TextCell "D5" "552.14"
Set-TextCell "E5" "  -1.30%  "

Set-TextCell "D6" "153.87"
Set-TextCell "E6" "  -2.64%  "

Set-TextCell "E7" "  +0.09%  "

Set-TextCell "D8" "0.595"
Set-TextCell "E8" "  +0.21%  "

Set-TextCell "E9" "  -3.07%  "

Set-TextCell "E10" "  -1.98%  "

Set-TextCell "D11" "5.43"
Set-TextCell "E11" "  -3.82%  "

Set-TextCell "E12" "  -2.92%  "

Set-TextCell "D13" "3.031.03"
Set-TextCell "E13" "  -5.33%  "

Set-TextCell "D14" "25.33"
Set-TextCell "E14" "  -4.69%  "

Set-TextCell "D15" "61.515.65"
Set-TextCell "E15" "  -2.97%  "

Set-TextCell "E16" "  -2.71%  "

Set-TextCell "D17" "2.579.70"
Set-TextCell "E17" "  -5.30%  "

Set-TextCell "D18" "11.53"
Set-TextCell "E18" "  -5.77%  "

Set-TextCell "D19" "4.52"
Set-TextCell "E19" "  -3.09%  "

Set-TextCell "D20" "336.54"
Set-TextCell "E20" "  -4.04%  "

Set-TextCell "D21" "6.03"
Set-TextCell "E21" "  -6.38%  "

Set-TextCell "D22" "0.998"
Set-TextCell "E22" "  -0.21%  "

Set-TextCell "D24" "62.92"
Set-TextCell "E24" "  -2.38%  "

Set-TextCell "E25" "  -0.99%  "

Set-TextCell "D26" "0.997"
Set-TextCell "E26" "  -0.20%  "

Set-TextCell "D27" "8.02"
Set-TextCell "E27" "  -2.31%  "

Set-TextCell "D28" "0.0₃0833"
Set-TextCell "E28" "  -5.70%  "

Set-TextCell "B29" "Aptos"
Set-TextCell "C29" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D29" "7.20"
Set-TextCell "E29" "  +0.74%  "

Set-TextCell "B30" "PancakeSwap"
Set-TextCell "C30" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D30" "1.91"
Set-TextCell "E30" "  -2.08%  "

Set-TextCell "E31" "  -4.79%  "

Set-TextCell "D33" "158.06"
Set-TextCell "E33" "  -4.40%  "

Set-TextCell "D34" "19.15"
Set-TextCell "E34" "  -3.66%  "

Set-TextCell "D35" "4.63"
Set-TextCell "E35" "  -3.93%  "

Set-TextCell "E36" "  -5.31%  "

Set-TextCell "E37" "  -1.05%  "

Set-TextCell "D38" "333.07"
Set-TextCell "E38" "  -3.34%  "

Set-TextCell "D39" "5.99"
Set-TextCell "E39" "  -1.39%  "

Set-TextCell "D40" "0.899"
Set-TextCell "E40" "  -6.43%  "

Set-TextCell "D41" "3.93"
Set-TextCell "E41" "  -2.51%  "

Set-TextCell "D42" "37.41"
Set-TextCell "E42" "  -2.12%  "

Set-TextCell "D43" "0.998"
Set-TextCell "E43" "  +0.05%  "

Set-TextCell "D44" "2.142.41"
Set-TextCell "E44" "  +0.99%  "

Set-TextCell "D45" "20.32"
Set-TextCell "E45" "  -4.98%  "

Set-TextCell "D46" "0.604"
Set-TextCell "E46" "  -3.57%  "

Set-TextCell "D47" "10.92"
Set-TextCell "E47" "  -1.09%  "

Set-TextCell "E48" "  -4.65%  "

Set-TextCell "D49" "19.53"
Set-TextCell "E49" "  -6.08%  "

Set-TextCell "D50" "0.0965"
Set-TextCell "E50" "  -1.90%  "

Set-TextCell "E51" "  -2.56%  "
